$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 3.366139666666667
$ws.Range("H2").Value = 10.098419
$ws.Range("I2").Value = 0.01725116351498256
$ws.Range("J2").Value = 0.01815407111703398
$ws.Range("M2").Value = 16.14072933333334
$ws.Range("N2").Value = 48.42218800000001
$ws.Range("O2").Value = 0.03423048004954622
$ws.Range("P2").Value = 0.03634868370049611
$ws.Range("Q2").Value = 54.33194925786357
$ws.Range("R2").Value = 488.987543320772
$ws.Range("S2").Value = 0.0005905156085310702
$ws.Range("T2").Value = 0.0006598765889093801
$ws.Range("G3").Value = 3.366139666666667
$ws.Range("H3").Value = 10.098419
$ws.Range("I3").Value = 0.01725116351498256
$ws.Range("J3").Value = 0.01815407111703398
$ws.Range("O3").Value = 0.1719151703242873
$ws.Range("P3").Value = 0.1825533892714798
$ws.Range("Q3").Value = 272.8704446211812
$ws.Range("R3").Value = 2455.834001590631
$ws.Range("S3").Value = 0.002965736713970358
$ws.Range("T3").Value = 0.003314087211490033
$ws.Range("G4").Value = 3.366139666666667
$ws.Range("H4").Value = 10.098419
$ws.Range("I4").Value = 0.01725116351498256
$ws.Range("J4").Value = 0.01815407111703398
$ws.Range("M4").Value = 168.70371
$ws.Range("N4").Value = 506.11113
$ws.Range("O4").Value = 0.3577786889414888
$ws.Range("P4").Value = 0.3799182594076638
$ws.Range("Q4").Value = 567.88025014483
$ws.Range("R4").Value = 5110.92225130347
$ws.Range("S4").Value = 0.006172098665105706
$ws.Range("T4").Value = 0.006897063099946492
$ws.Range("G5").Value = 3.366139666666667
$ws.Range("H5").Value = 10.098419
$ws.Range("I5").Value = 0.01725116351498256
$ws.Range("J5").Value = 0.01815407111703398
$ws.Range("M5").Value = 82.43477250000001
$ws.Range("N5").Value = 164.869545
$ws.Range("O5").Value = 0.1748236883957081
$ws.Range("P5").Value = 0.1237612588479007
$ws.Range("Q5").Value = 277.4869576248925
$ws.Range("R5").Value = 1664.921745749355
$ws.Range("S5").Value = 0.003015912034806719
$ws.Range("T5").Value = 0.00224677069465844
$ws.Range("G6").Value = 3.366139666666667
$ws.Range("H6").Value = 10.098419
$ws.Range("I6").Value = 0.01725116351498256
$ws.Range("J6").Value = 0.01815407111703398
$ws.Range("M6").Value = 123.1883796666667
$ws.Range("N6").Value = 369.565139
$ws.Range("O6").Value = 0.2612519722889696
$ws.Range("P6").Value = 0.2774184087724594
$ws.Range("Q6").Value = 414.6692912683601
$ws.Range("R6").Value = 3732.023621415241
$ws.Range("S6").Value = 0.004506900492568707
$ws.Range("T6").Value = 0.005036273522029631
$ws.Range("I7").Value = 0.7504462978934635
$ws.Range("J7").Value = 0.7897238612132288
$ws.Range("M7").Value = 16.14072933333334
$ws.Range("N7").Value = 48.42218800000001
$ws.Range("O7").Value = 0.03423048004954622
$ws.Range("P7").Value = 0.03634868370049611
$ws.Range("Q7").Value = 2363.504939390776
$ws.Range("R7").Value = 21271.54445451698
$ws.Range("S7").Value = 0.02568813702829802
$ws.Range("T7").Value = 0.02870542284197414
$ws.Range("I8").Value = 0.7504462978934635
$ws.Range("J8").Value = 0.7897238612132288
$ws.Range("O8").Value = 0.1719151703242873
$ws.Range("P8").Value = 0.1825533892714798
$ws.Range("S8").Value = 0.1290131031215856
$ws.Range("T8").Value = 0.1441667674530347
$ws.Range("I9").Value = 0.7504462978934635
$ws.Range("J9").Value = 0.7897238612132288
$ws.Range("M9").Value = 168.70371
$ws.Range("N9").Value = 506.11113
$ws.Range("O9").Value = 0.3577786889414888
$ws.Range("P9").Value = 0.3799182594076638
$ws.Range("Q9").Value = 24703.47179759095
$ws.Range("R9").Value = 222331.2461783186
$ws.Range("S9").Value = 0.2684936925813173
$ws.Range("T9").Value = 0.3000305147648293
$ws.Range("I10").Value = 0.7504462978934635
$ws.Range("J10").Value = 0.7897238612132288
$ws.Range("M10").Value = 82.43477250000001
$ws.Range("N10").Value = 164.869545
$ws.Range("O10").Value = 0.1748236883957081
$ws.Range("P10").Value = 0.1237612588479007
$ws.Range("Q10").Value = 12071.01537716377
$ws.Range("R10").Value = 72426.09226298259
$ws.Range("S10").Value = 0.1311957897406396
$ws.Range("T10").Value = 0.09773721920597402
$ws.Range("I11").Value = 0.7504462978934635
$ws.Range("J11").Value = 0.7897238612132288
$ws.Range("M11").Value = 123.1883796666667
$ws.Range("N11").Value = 369.565139
$ws.Range("O11").Value = 0.2612519722889696
$ws.Range("P11").Value = 0.2774184087724594
$ws.Range("Q11").Value = 18038.6113790054
$ws.Range("R11").Value = 162347.5024110486
$ws.Range("S11").Value = 0.1960555754216229
$ws.Range("T11").Value = 0.2190839369474165
$ws.Range("G12").Value = 14.89209833333333
$ws.Range("H12").Value = 44.676295
$ws.Range("I12").Value = 0.07632066665966204
$ws.Range("J12").Value = 0.08031520940808551
$ws.Range("M12").Value = 16.14072933333334
$ws.Range("N12").Value = 48.42218800000001
$ws.Range("O12").Value = 0.03423048004954622
$ws.Range("P12").Value = 0.03634868370049611
$ws.Range("Q12").Value = 240.3693284037178
$ws.Range("R12").Value = 2163.32395563346
$ws.Range("S12").Value = 0.002612493057461629
$ws.Range("T12").Value = 0.00291935214311361
$ws.Range("G13").Value = 14.89209833333333
$ws.Range("H13").Value = 44.676295
$ws.Range("I13").Value = 0.07632066665966204
$ws.Range("J13").Value = 0.08031520940808551
$ws.Range("O13").Value = 0.1719151703242873
$ws.Range("P13").Value = 0.1825533892714798
$ws.Range("Q13").Value = 1207.202878062106
$ws.Range("R13").Value = 10864.82590255895
$ws.Range("S13").Value = 0.01312068040805896
$ws.Range("T13").Value = 0.01466181368749465
$ws.Range("G14").Value = 14.89209833333333
$ws.Range("H14").Value = 44.676295
$ws.Range("I14").Value = 0.07632066665966204
$ws.Range("J14").Value = 0.08031520940808551
$ws.Range("M14").Value = 168.70371
$ws.Range("N14").Value = 506.11113
$ws.Range("O14").Value = 0.3577786889414888
$ws.Range("P14").Value = 0.3799182594076638
$ws.Range("Q14").Value = 2512.35223851815
$ws.Range("R14").Value = 22611.17014666335
$ws.Range("S14").Value = 0.02730590805663428
$ws.Range("T14").Value = 0.03051321456228187
$ws.Range("G15").Value = 14.89209833333333
$ws.Range("H15").Value = 44.676295
$ws.Range("I15").Value = 0.07632066665966204
$ws.Range("J15").Value = 0.08031520940808551
$ws.Range("M15").Value = 82.43477250000001
$ws.Range("N15").Value = 164.869545
$ws.Range("O15").Value = 0.1748236883957081
$ws.Range("P15").Value = 0.1237612588479007
$ws.Range("Q15").Value = 1227.626738155962
$ws.Range("R15").Value = 7365.760428935775
$ws.Range("S15").Value = 0.01334266044626146
$ws.Range("T15").Value = 0.009939911420977419
$ws.Range("G16").Value = 14.89209833333333
$ws.Range("H16").Value = 44.676295
$ws.Range("I16").Value = 0.07632066665966204
$ws.Range("J16").Value = 0.08031520940808551
$ws.Range("M16").Value = 123.1883796666667
$ws.Range("N16").Value = 369.565139
$ws.Range("O16").Value = 0.2612519722889696
$ws.Range("P16").Value = 0.2774184087724594
$ws.Range("Q16").Value = 1834.53346352
$ws.Range("R16").Value = 16510.80117168
$ws.Range("S16").Value = 0.01993892469124571
$ws.Range("T16").Value = 0.02228091759421795
$ws.Range("G17").Value = 29.1141605
$ws.Range("H17").Value = 58.22832099999999
$ws.Range("I17").Value = 0.1492074581338761
$ws.Range("J17").Value = 0.1046778788302885
$ws.Range("M17").Value = 16.14072933333334
$ws.Range("N17").Value = 48.42218800000001
$ws.Range("O17").Value = 0.03423048004954622
$ws.Range("P17").Value = 0.03634868370049611
$ws.Range("Q17").Value = 469.9237843977247
$ws.Range("R17").Value = 2819.542706386348
$ws.Range("S17").Value = 0.005107442918895149
$ws.Range("T17").Value = 0.003804903108041014
$ws.Range("G18").Value = 29.1141605
$ws.Range("H18").Value = 58.22832099999999
$ws.Range("I18").Value = 0.1492074581338761
$ws.Range("J18").Value = 0.1046778788302885
$ws.Range("O18").Value = 0.1719151703242873
$ws.Range("P18").Value = 0.1825533892714798
$ws.Range("Q18").Value = 2360.090402390938
$ws.Range("R18").Value = 14160.54241434563
$ws.Range("S18").Value = 0.02565102557873928
$ws.Range("T18").Value = 0.01910930156221845
$ws.Range("G19").Value = 29.1141605
$ws.Range("H19").Value = 58.22832099999999
$ws.Range("I19").Value = 0.1492074581338761
$ws.Range("J19").Value = 0.1046778788302885
$ws.Range("M19").Value = 168.70371
$ws.Range("N19").Value = 506.11113
$ws.Range("O19").Value = 0.3577786889414888
$ws.Range("P19").Value = 0.3799182594076638
$ws.Range("Q19").Value = 4911.666889885454
$ws.Range("R19").Value = 29470.00133931273
$ws.Range("S19").Value = 0.05338324875143027
$ws.Range("T19").Value = 0.03976903752368954
$ws.Range("G20").Value = 29.1141605
$ws.Range("H20").Value = 58.22832099999999
$ws.Range("I20").Value = 0.1492074581338761
$ws.Range("J20").Value = 0.1046778788302885
$ws.Range("M20").Value = 82.43477250000001
$ws.Range("N20").Value = 164.869545
$ws.Range("O20").Value = 0.1748236883957081
$ws.Range("P20").Value = 0.1237612588479007
$ws.Range("Q20").Value = 2400.019197345986
$ws.Range("R20").Value = 9600.076789383946
$ws.Range("S20").Value = 0.02608499816711241
$ws.Range("T20").Value = 0.01295506605756452
$ws.Range("G21").Value = 29.1141605
$ws.Range("H21").Value = 58.22832099999999
$ws.Range("I21").Value = 0.1492074581338761
$ws.Range("J21").Value = 0.1046778788302885
$ws.Range("M21").Value = 123.1883796666667
$ws.Range("N21").Value = 369.565139
$ws.Range("O21").Value = 0.2612519722889696
$ws.Range("P21").Value = 0.2774184087724594
$ws.Range("Q21").Value = 3586.526257350269
$ws.Range("R21").Value = 21519.15754410162
$ws.Range("S21").Value = 0.03898074271769899
$ws.Range("T21").Value = 0.02903957057877495
$ws.Range("G22").Value = 1.32186
$ws.Range("H22").Value = 3.96558
$ws.Range("I22").Value = 0.006774413798015763
$ws.Range("J22").Value = 0.007128979431363227
$ws.Range("M22").Value = 16.14072933333334
$ws.Range("N22").Value = 48.42218800000001
$ws.Range("O22").Value = 0.03423048004954622
$ws.Range("P22").Value = 0.03634868370049611
$ws.Range("Q22").Value = 21.33578447656
$ws.Range("R22").Value = 192.02206028904
$ws.Range("S22").Value = 0.0002318914363603492
$ws.Range("T22").Value = 0.0002591290184579646
$ws.Range("G23").Value = 1.32186
$ws.Range("H23").Value = 3.96558
$ws.Range("I23").Value = 0.006774413798015763
$ws.Range("J23").Value = 0.007128979431363227
$ws.Range("O23").Value = 0.1719151703242873
$ws.Range("P23").Value = 0.1825533892714798
$ws.Range("Q23").Value = 107.15435532838
$ws.Range("R23").Value = 964.3891979554201
$ws.Range("S23").Value = 0.001164624501933082
$ws.Range("T23").Value = 0.001301419357242024
$ws.Range("G24").Value = 1.32186
$ws.Range("H24").Value = 3.96558
$ws.Range("I24").Value = 0.006774413798015763
$ws.Range("J24").Value = 0.007128979431363227
$ws.Range("M24").Value = 168.70371
$ws.Range("N24").Value = 506.11113
$ws.Range("O24").Value = 0.3577786889414888
$ws.Range("P24").Value = 0.3799182594076638
$ws.Range("Q24").Value = 223.0026861006
$ws.Range("R24").Value = 2007.0241749054
$ws.Range("S24").Value = 0.002423740887001211
$ws.Range("T24").Value = 0.002708429456916554
$ws.Range("G25").Value = 1.32186
$ws.Range("H25").Value = 3.96558
$ws.Range("I25").Value = 0.006774413798015763
$ws.Range("J25").Value = 0.007128979431363227
$ws.Range("M25").Value = 82.43477250000001
$ws.Range("N25").Value = 164.869545
$ws.Range("O25").Value = 0.1748236883957081
$ws.Range("P25").Value = 0.1237612588479007
$ws.Range("Q25").Value = 108.96722837685
$ws.Range("R25").Value = 653.8033702611001
$ws.Range("S25").Value = 0.001184328006887893
$ws.Range("T25").Value = 0.0008822914687263042
$ws.Range("G26").Value = 1.32186
$ws.Range("H26").Value = 3.96558
$ws.Range("I26").Value = 0.006774413798015763
$ws.Range("J26").Value = 0.007128979431363227
$ws.Range("M26").Value = 123.1883796666667
$ws.Range("N26").Value = 369.565139
$ws.Range("O26").Value = 0.2612519722889696
$ws.Range("P26").Value = 0.2774184087724594
$ws.Range("Q26").Value = 162.83779154618
$ws.Range("R26").Value = 1465.54012391562
$ws.Range("S26").Value = 0.001769828965833227
$ws.Range("T26").Value = 0.001977710130020379
